$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 16.59582222408542
$ws.Range("C2").Value = 10.70525454909819
$ws.Range("D2").Value = 9.898593645394611
$ws.Range("F2").Value = 30.64735522567415
$ws.Range("G2").Value = 30.52755002881299
$ws.Range("H2").Value = 14.72558321375142
$ws.Range("J2").Value = 10.45015723247472
$ws.Range("L2").Value = 11.57082506283825
$ws.Range("N2").Value = 17.64349507048993
$ws.Range("O2").Value = 22.65447234228262
$ws.Range("B3").Value = 16.1434667583125
$ws.Range("C3").Value = 10.59953461193187
$ws.Range("D3").Value = 9.900526618958461
$ws.Range("F3").Value = 30.70071557058448
$ws.Range("G3").Value = 30.54584793298028
$ws.Range("H3").Value = 14.76839332506794
$ws.Range("J3").Value = 10.47708063343948
$ws.Range("L3").Value = 11.55240839342511
$ws.Range("N3").Value = 17.68371475885496
$ws.Range("O3").Value = 22.71400486308479
$ws.Range("B4").Value = 15.86098299596559
$ws.Range("C4").Value = 10.53399766320992
$ws.Range("D4").Value = 9.903018563756074
$ws.Range("F4").Value = 30.74077767560525
$ws.Range("G4").Value = 30.56695297631545
$ws.Range("H4").Value = 14.79712866025968
$ws.Range("J4").Value = 10.49455252309765
$ws.Range("L4").Value = 11.54251427557426
$ws.Range("N4").Value = 17.71017479455948
$ws.Range("O4").Value = 22.75558660528883
$ws.Range("B5").Value = 15.74484503987587
$ws.Range("C5").Value = 10.50715134606009
$ws.Range("D5").Value = 9.904363164122614
$ws.Range("F5").Value = 30.75893554952356
$ws.Range("G5").Value = 30.57803063397859
$ws.Range("H5").Value = 14.80945426412243
$ws.Range("J5").Value = 10.50190959232477
$ws.Range("L5").Value = 11.53884083322959
$ws.Range("N5").Value = 17.72140202363628
$ws.Range("O5").Value = 22.77379342729816
$ws.Range("B6").Value = 15.72550368770051
$ws.Range("C6").Value = 10.50268562596113
$ws.Range("D6").Value = 9.90460633722412
$ws.Range("F6").Value = 30.76206121280251
$ws.Range("G6").Value = 30.58001951378466
$ws.Range("H6").Value = 14.81153810051418
$ws.Range("J6").Value = 10.50314556671115
$ws.Range("L6").Value = 11.53825259614706
$ws.Range("N6").Value = 17.72329316784857
$ws.Range("O6").Value = 22.7768927951371
$ws.Range("B7").Value = 15.85942063467497
$ws.Range("C7").Value = 10.53363614588352
$ws.Range("D7").Value = 9.903035363763028
$ws.Range("F7").Value = 30.74101514484735
$ws.Range("G7").Value = 30.56709235102965
$ws.Range("H7").Value = 14.79729239503171
$ws.Range("J7").Value = 10.49465078222284
$ws.Range("L7").Value = 11.54246327893469
$ws.Range("N7").Value = 17.7103244079578
$ws.Range("O7").Value = 22.75582704289029
$ws.Range("B8").Value = 16.44093015362169
$ws.Range("C8").Value = 10.66894241847883
$ws.Range("D8").Value = 9.898989849343245
$ws.Range("F8").Value = 30.66423736656704
$ws.Range("G8").Value = 30.53180830035782
$ws.Range("H8").Value = 14.73983552637769
$ws.Range("J8").Value = 10.45924549133245
$ws.Range("L8").Value = 11.5641834070039
$ws.Range("N8").Value = 17.65699693301576
$ws.Range("O8").Value = 22.67395401401322
$ws.Range("B9").Value = 17.53687365772846
$ws.Range("C9").Value = 10.92856050696507
$ws.Range("D9").Value = 9.901368102858571
$ws.Range("F9").Value = 30.57170211399112
$ws.Range("G9").Value = 30.54106288904535
$ws.Range("H9").Value = 14.64661292257668
$ws.Range("J9").Value = 10.39725532328552
$ws.Range("L9").Value = 11.61785633679192
$ws.Range("N9").Value = 17.56639403720446
$ws.Range("O9").Value = 22.55340662255061
$ws.Range("B10").Value = 18.30664254789727
$ws.Range("C10").Value = 11.11481214567931
$ws.Range("D10").Value = 9.90933879147315
$ws.Range("F10").Value = 30.53921297711794
$ws.Range("G10").Value = 30.59574892915078
$ws.Range("H10").Value = 14.58999585269405
$ws.Range("J10").Value = 10.35621088959003
$ws.Range("L10").Value = 11.66384933621456
$ws.Range("N10").Value = 17.50830226268312
$ws.Range("O10").Value = 22.48936273297079
$ws.Range("B11").Value = 18.64759683735362
$ws.Range("C11").Value = 11.19835991752598
$ws.Range("D11").Value = 9.914300964148751
$ws.Range("F11").Value = 30.53215355332251
$ws.Range("G11").Value = 30.6309986717108
$ws.Range("H11").Value = 14.56682061595121
$ws.Range("J11").Value = 10.33850826947089
$ws.Range("L11").Value = 11.68615161047819
$ws.Range("N11").Value = 17.4837059060187
$ws.Range("O11").Value = 22.46557625974493
$ws.Range("B12").Value = 18.77526837661908
$ws.Range("C12").Value = 11.22981066644165
$ws.Range("D12").Value = 9.916370688695268
$ws.Range("F12").Value = 30.53059030880463
$ws.Range("G12").Value = 30.64583289041795
$ws.Range("H12").Value = 14.55841601666034
$ws.Range("J12").Value = 10.33194348647029
$ws.Range("L12").Value = 11.69479113203996
$ws.Range("N12").Value = 17.4746543734216
$ws.Range("O12").Value = 22.45733942396168
$ws.Range("B13").Value = 18.74783780707149
$ws.Range("C13").Value = 11.22304577403708
$ws.Range("D13").Value = 9.915916481821151
$ws.Range("F13").Value = 30.53087762336305
$ws.Range("G13").Value = 30.64257210188785
$ws.Range("H13").Value = 14.56020957611143
$ws.Range("J13").Value = 10.33335116373154
$ws.Range("L13").Value = 11.69292189244637
$ws.Range("N13").Value = 17.47659211366009
$ws.Range("O13").Value = 22.45907907885688
$ws.Range("B14").Value = 18.6581299407963
$ws.Range("C14").Value = 11.20095123528906
$ws.Range("D14").Value = 9.914467429494787
$ws.Range("F14").Value = 30.53200270137022
$ws.Range("G14").Value = 30.63218935176769
$ws.Range("H14").Value = 14.56612172007807
$ws.Range("J14").Value = 10.33796540150458
$ws.Range("L14").Value = 11.68685852355483
$ws.Range("N14").Value = 17.48295597177658
$ws.Range("O14").Value = 22.4648831582518
$ws.Range("B15").Value = 18.60299046497721
$ws.Range("C15").Value = 11.18739283178629
$ws.Range("D15").Value = 9.913604626237687
$ws.Range("F15").Value = 30.53283638362353
$ws.Range("G15").Value = 30.62602292322359
$ws.Range("H15").Value = 14.5697914504886
$ws.Range("J15").Value = 10.34080981897312
$ws.Range("L15").Value = 11.68316968682541
$ws.Range("N15").Value = 17.48688819708267
$ws.Range("O15").Value = 22.46853872121933
$ws.Range("B16").Value = 18.28416482673488
$ws.Range("C16").Value = 11.10932687735689
$ws.Range("D16").Value = 9.909041265303543
$ws.Range("F16").Value = 30.53982967967399
$ws.Range("G16").Value = 30.59365356232755
$ws.Range("H16").Value = 14.59156235766802
$ws.Range("J16").Value = 10.35738724873928
$ws.Range("L16").Value = 11.66241923384862
$ws.Range("N16").Value = 17.50994646607257
$ws.Range("O16").Value = 22.49102498449839
$ws.Range("B17").Value = 18.08612953665052
$ws.Range("C17").Value = 11.06112143058745
$ws.Range("D17").Value = 9.906583008240528
$ws.Range("F17").Value = 30.54609726663654
$ws.Range("G17").Value = 30.57644932006324
$ws.Range("H17").Value = 14.60557919210038
$ws.Range("J17").Value = 10.36780472651003
$ws.Range("L17").Value = 11.65003992105763
$ws.Range("N17").Value = 17.5245602404342
$ws.Range("O17").Value = 22.50619042882429
$ws.Range("B18").Value = 17.97136345760324
$ws.Range("C18").Value = 11.03328497492388
$ws.Range("D18").Value = 9.905294955360333
$ws.Range("F18").Value = 30.5504289170824
$ws.Range("G18").Value = 30.56753069928845
$ws.Range("H18").Value = 14.61388414230387
$ws.Range("J18").Value = 10.37388779115013
$ws.Range("L18").Value = 11.64304981186066
$ws.Range("N18").Value = 17.53313796679237
$ws.Range("O18").Value = 22.51541643449357
$ws.Range("B19").Value = 17.93236136805805
$ws.Range("C19").Value = 11.02384168191557
$ws.Range("D19").Value = 9.904880505297587
$ws.Range("F19").Value = 30.55202034310127
$ws.Range("G19").Value = 30.56467891498154
$ws.Range("H19").Value = 14.61673775262523
$ws.Range("J19").Value = 10.37596309200899
$ws.Range("L19").Value = 11.64070556086572
$ws.Range("N19").Value = 17.53607184240263
$ws.Range("O19").Value = 22.51862657974385
$ws.Range("B20").Value = 18.10730075226402
$ws.Range("C20").Value = 11.06626449049704
$ws.Range("D20").Value = 9.906831677194337
$ws.Range("F20").Value = 30.5453548591946
$ws.Range("G20").Value = 30.57817968023889
$ws.Range("H20").Value = 14.60406194080669
$ws.Range("J20").Value = 10.3666863321555
$ws.Range("L20").Value = 11.6513442810609
$ws.Range("N20").Value = 17.52298675402637
$ws.Range("O20").Value = 22.50452394488624
$ws.Range("B21").Value = 18.6845192719975
$ws.Range("C21").Value = 11.20744614833596
$ws.Range("D21").Value = 9.914887889316216
$ws.Range("F21").Value = 30.53164211802849
$ws.Range("G21").Value = 30.63519874539331
$ws.Range("H21").Value = 14.56437509894157
$ws.Range("J21").Value = 10.33660632443498
$ws.Range("L21").Value = 11.68863424826285
$ws.Range("N21").Value = 17.48107963060777
$ws.Range("O21").Value = 22.46315743432205
$ws.Range("B22").Value = 19.05332307745502
$ws.Range("C22").Value = 11.29861888212042
$ws.Range("D22").Value = 9.921263623171138
$ws.Range("F22").Value = 30.52914944365796
$ws.Range("G22").Value = 30.68112094321479
$ws.Range("H22").Value = 14.54060215840078
$ws.Range("J22").Value = 10.31775620131318
$ws.Range("L22").Value = 11.71413484750691
$ws.Range("N22").Value = 17.45522113390881
$ws.Range("O22").Value = 22.44061415084669
$ws.Range("B23").Value = 18.85729319825918
$ws.Range("C23").Value = 11.25006436260347
$ws.Range("D23").Value = 9.917759679999213
$ws.Range("F23").Value = 30.52988810962419
$ws.Range("G23").Value = 30.65582168946729
$ws.Range("H23").Value = 14.5530920498919
$ws.Range("J23").Value = 10.32774301014003
$ws.Range("L23").Value = 11.70042281459074
$ws.Range("N23").Value = 17.46888246574311
$ws.Range("O23").Value = 22.45223442800472
$ws.Range("B24").Value = 18.09773208046349
$ws.Range("C24").Value = 11.06393969336096
$ws.Range("D24").Value = 9.906718863774465
$ws.Range("F24").Value = 30.54568823289972
$ws.Range("G24").Value = 30.57739435526984
$ws.Range("H24").Value = 14.60474712234929
$ws.Range("J24").Value = 10.36719166566049
$ws.Range("L24").Value = 11.65075418423534
$ws.Range("N24").Value = 17.5236975787293
$ws.Range("O24").Value = 22.50527578253326
$ws.Range("B25").Value = 17.2460512070159
$ws.Range("C25").Value = 10.85904361075941
$ws.Range("D25").Value = 9.899626903219394
$ws.Range("F25").Value = 30.59050808282389
$ws.Range("G25").Value = 30.53014687111805
$ws.Range("H25").Value = 14.66974806605624
$ws.Range("J25").Value = 10.41323255718859
$ws.Range("L25").Value = 11.60216956293792
$ws.Range("N25").Value = 17.58941327309739
$ws.Range("O25").Value = 22.58172027420047
